# Update basement (type_base) construction type for STANDARD2..STANDARD6
# from FLOOR_AS4 to FLOOR_AS6 on the ENVELOPE_ASSEMBLIES sheet (column I,
# rows 3-7). STANDARD1 (row 2) keeps its existing FLOOR_AS4 value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ENVELOPE_ASSEMBLIES")

$ws.Range("I3:I7").Value = "FLOOR_AS6"

# Reflect the updated active selection on this sheet, matching the
# author's final cursor position.
$ws.Activate()
$ws.Range("H17").Select()
